$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 update
$ws.Range("G2").Value = 2.1

# Row 5 updates
$ws.Range("G5").Value = 1.91
$ws.Range("I5").Value = 4.1
$ws.Range("L5").Value = 1.42
$ws.Range("M5").Value = 2.47
$ws.Range("N5").Value = 2.2
$ws.Range("O5").Value = 1.53
$ws.Range("P5").Value = 1.47
$ws.Range("Q5").Value = 2.35
$ws.Range("R5").Value = 1.98
$ws.Range("S5").Value = 1.65
$ws.Range("T5").Value = 5.7
$ws.Range("U5").Value = 8
$ws.Range("W5").Value = 16.5
$ws.Range("Y5").Value = 35
$ws.Range("Z5").Value = 7.1
$ws.Range("AA5").Value = 6.1
$ws.Range("AB5").Value = 17.5
$ws.Range("AC5").Value = 110
$ws.Range("AE5").Value = 9.25
$ws.Range("AF5").Value = 21
$ws.Range("AG5").Value = 14
$ws.Range("AH5").Value = 70
$ws.Range("AI5").Value = 45
$ws.Range("AJ5").Value = 60
